$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows at 658, pushing the existing 658-723 block down to 667-732
$ws.Rows("658:666").Insert()

# Fill the newly inserted rows with the new match data
$ws.Range("A658").Value2 = 540458
$ws.Range("B658").Value2 = "2025-10-03T18:30:00Z"
$ws.Range("C658").Value2 = 6
$ws.Range("D658").Value2 = 2
$ws.Range("E658").Value2 = "TSG 1899 Hoffenheim"
$ws.Range("F658").Value2 = 1
$ws.Range("G658").Value2 = "1. FC Köln"
$ws.Range("H658").Value2 = 0
$ws.Range("I658").Value2 = 1
$ws.Range("J658").Value2 = "AwayWin"
$ws.Range("K658").Value2 = 3
$ws.Range("L658").Value2 = 4

$ws.Range("A659").Value2 = 540451
$ws.Range("B659").Value2 = "2025-10-04T13:30:00Z"
$ws.Range("C659").Value2 = 6
$ws.Range("D659").Value2 = 3
$ws.Range("E659").Value2 = "Bayer 04 Leverkusen"
$ws.Range("F659").Value2 = 28
$ws.Range("G659").Value2 = "1. FC Union Berlin"
$ws.Range("H659").Value2 = 2
$ws.Range("I659").Value2 = 0
$ws.Range("J659").Value2 = "HomeWin"
$ws.Range("K659").Value2 = 1
$ws.Range("L659").Value2 = 3

$ws.Range("A660").Value2 = 540453
$ws.Range("B660").Value2 = "2025-10-04T13:30:00Z"
$ws.Range("C660").Value2 = 6
$ws.Range("D660").Value2 = 4
$ws.Range("E660").Value2 = "Borussia Dortmund"
$ws.Range("F660").Value2 = 721
$ws.Range("G660").Value2 = "RB Leipzig"
$ws.Range("H660").Value2 = 1
$ws.Range("I660").Value2 = 1
$ws.Range("J660").Value2 = "Draw"
$ws.Range("K660").Value2 = 6
$ws.Range("L660").Value2 = 7

$ws.Range("A661").Value2 = 540454
$ws.Range("B661").Value2 = "2025-10-04T13:30:00Z"
$ws.Range("C661").Value2 = 6
$ws.Range("D661").Value2 = 12
$ws.Range("E661").Value2 = "SV Werder Bremen"
$ws.Range("F661").Value2 = 20
$ws.Range("G661").Value2 = "FC St. Pauli"
$ws.Range("H661").Value2 = 1
$ws.Range("I661").Value2 = 0
$ws.Range("J661").Value2 = "HomeWin"
$ws.Range("K661").Value2 = 3
$ws.Range("L661").Value2 = 3

$ws.Range("A662").Value2 = 540457
$ws.Range("B662").Value2 = "2025-10-04T13:30:00Z"
$ws.Range("C662").Value2 = 6
$ws.Range("D662").Value2 = 16
$ws.Range("E662").Value2 = "FC Augsburg"
$ws.Range("F662").Value2 = 11
$ws.Range("G662").Value2 = "VfL Wolfsburg"
$ws.Range("H662").Value2 = 3
$ws.Range("I662").Value2 = 1
$ws.Range("J662").Value2 = "HomeWin"
$ws.Range("K662").Value2 = 0
$ws.Range("L662").Value2 = 0

$ws.Range("A663").Value2 = 540452
$ws.Range("B663").Value2 = "2025-10-04T16:30:00Z"
$ws.Range("C663").Value2 = 6
$ws.Range("D663").Value2 = 19
$ws.Range("E663").Value2 = "Eintracht Frankfurt"
$ws.Range("F663").Value2 = 5
$ws.Range("G663").Value2 = "FC Bayern München"
$ws.Range("H663").Value2 = 0
$ws.Range("I663").Value2 = 3
$ws.Range("J663").Value2 = "AwayWin"
$ws.Range("K663").Value2 = 0
$ws.Range("L663").Value2 = 6

$ws.Range("A664").Value2 = 540455
$ws.Range("B664").Value2 = "2025-10-05T13:30:00Z"
$ws.Range("C664").Value2 = 6
$ws.Range("D664").Value2 = 10
$ws.Range("E664").Value2 = "VfB Stuttgart"
$ws.Range("F664").Value2 = 44
$ws.Range("G664").Value2 = "1. FC Heidenheim 1846"
$ws.Range("H664").Value2 = 1
$ws.Range("I664").Value2 = 0
$ws.Range("J664").Value2 = "HomeWin"
$ws.Range("K664").Value2 = 3
$ws.Range("L664").Value2 = 0

$ws.Range("A665").Value2 = 540459
$ws.Range("B665").Value2 = "2025-10-05T15:30:00Z"
$ws.Range("C665").Value2 = 6
$ws.Range("D665").Value2 = 7
$ws.Range("E665").Value2 = "Hamburger SV"
$ws.Range("F665").Value2 = 15
$ws.Range("G665").Value2 = "1. FSV Mainz 05"
$ws.Range("H665").Value2 = 4
$ws.Range("I665").Value2 = 0
$ws.Range("J665").Value2 = "HomeWin"
$ws.Range("K665").Value2 = 6
$ws.Range("L665").Value2 = 1

$ws.Range("A666").Value2 = 540456
$ws.Range("B666").Value2 = "2025-10-05T17:30:00Z"
$ws.Range("C666").Value2 = 6
$ws.Range("D666").Value2 = 18
$ws.Range("E666").Value2 = "Borussia Mönchengladbach"
$ws.Range("F666").Value2 = 17
$ws.Range("G666").Value2 = "SC Freiburg"
$ws.Range("H666").Value2 = 0
$ws.Range("I666").Value2 = 0
$ws.Range("J666").Value2 = "Draw"
$ws.Range("K666").Value2 = 1
$ws.Range("L666").Value2 = 1

Write-Output "done"
